$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 45962
$ws.Range("D8").Value = 167.82
$ws.Range("E8").Value = 161.06
$ws.Range("F8").Value = 171.06
$ws.Range("G8").Value = 161.22

# Row 9
$ws.Range("A9").Value = 45962
$ws.Range("D9").Value = 167.82
$ws.Range("E9").Value = 161.06
$ws.Range("F9").Value = 171.06
$ws.Range("G9").Value = 161.22

# Row 10
$ws.Range("A10").Value = 45962
$ws.Range("D10").Value = 169.92
$ws.Range("E10").Value = 163.41
$ws.Range("F10").Value = 173.41
$ws.Range("G10").Value = 163.9

# Row 11
$ws.Range("A11").Value = 45961
$ws.Range("D11").Value = 167.32
$ws.Range("E11").Value = 160.94
$ws.Range("F11").Value = 170.94
$ws.Range("G11").Value = 161.09

# Row 12
$ws.Range("A12").Value = 45961
$ws.Range("D12").Value = 167.32
$ws.Range("E12").Value = 160.94
$ws.Range("F12").Value = 170.94
$ws.Range("G12").Value = 161.09

# Row 13
$ws.Range("A13").Value = 45961
$ws.Range("D13").Value = 169.79
$ws.Range("E13").Value = 163.41999999999999
$ws.Range("F13").Value = 173.42
$ws.Range("G13").Value = 163.89

# Row 17
$ws.Range("A17").Value = 45962
$ws.Range("D17").Value = 173.43
$ws.Range("E17").Value = 166.39
$ws.Range("F17").Value = 176.39

# Row 18
$ws.Range("A18").Value = 45961
$ws.Range("D18").Value = 173.29
$ws.Range("E18").Value = 166.35
$ws.Range("F18").Value = 176.35

# Row 22
$ws.Range("A22").Value = 45962
$ws.Range("D22").Value = 169.06
$ws.Range("E22").Value = 162.43
$ws.Range("F22").Value = 172.03
$ws.Range("G22").Value = 163.72

# Row 23
$ws.Range("A23").Value = 45962
$ws.Range("D23").Value = 174.69
$ws.Range("E23").Value = 167.14
$ws.Range("F23").Value = 177.14

# Row 24
$ws.Range("A24").Value = 45962
$ws.Range("D24").Value = 174.51
$ws.Range("E24").Value = 167.31
$ws.Range("F24").Value = 177.31

# Row 25
$ws.Range("A25").Value = 45962
$ws.Range("D25").Value = 175.34
$ws.Range("E25").Value = 166.7
$ws.Range("F25").Value = 176.7
$ws.Range("G25").Value = 166.74

# Row 26
$ws.Range("A26").Value = 45962
$ws.Range("D26").Value = 174.08
$ws.Range("E26").Value = 168.26
$ws.Range("F26").Value = 178.26

# Row 27
$ws.Range("A27").Value = 45961
$ws.Range("D27").Value = 169
$ws.Range("E27").Value = 162.44
$ws.Range("F27").Value = 172.04
$ws.Range("G27").Value = 163.62

# Row 28
$ws.Range("A28").Value = 45961
$ws.Range("D28").Value = 174.56
$ws.Range("E28").Value = 167.15
$ws.Range("F28").Value = 177.15

# Row 29
$ws.Range("A29").Value = 45961
$ws.Range("D29").Value = 174.37
$ws.Range("E29").Value = 167.32
$ws.Range("F29").Value = 177.32

# Row 30
$ws.Range("A30").Value = 45961
$ws.Range("D30").Value = 175.2
$ws.Range("E30").Value = 166.71
$ws.Range("F30").Value = 176.71
$ws.Range("G30").Value = 166.54

# Row 31
$ws.Range("A31").Value = 45961
$ws.Range("D31").Value = 173.94
$ws.Range("E31").Value = 168.28
$ws.Range("F31").Value = 178.28

# Row 35
$ws.Range("A35").Value = 45962
$ws.Range("D35").Value = 168.29

# Row 36
$ws.Range("A36").Value = 45961
$ws.Range("D36").Value = 168.04
$ws.Range("E36").Value = 160.63
$ws.Range("F36").Value = 169.63

# Row 40
$ws.Range("A40").Value = 45962
$ws.Range("D40").Value = 173.77
$ws.Range("E40").Value = 166.01
$ws.Range("F40").Value = 176.01

# Row 41
$ws.Range("A41").Value = 45962
$ws.Range("D41").Value = 173.48
$ws.Range("E41").Value = 166.43
$ws.Range("F41").Value = 176.43

# Row 42
$ws.Range("A42").Value = 45961
$ws.Range("D42").Value = 173.64
$ws.Range("E42").Value = 165.98
$ws.Range("F42").Value = 175.98

# Row 43
$ws.Range("A43").Value = 45961
$ws.Range("D43").Value = 173.35
$ws.Range("E43").Value = 166.4
$ws.Range("F43").Value = 176.4

# Row 47
$ws.Range("A47").Value = 45962
$ws.Range("D47").Value = 167.91
$ws.Range("E47").Value = 161.74
$ws.Range("F47").Value = 171.74

# Row 48
$ws.Range("A48").Value = 45962
$ws.Range("D48").Value = 167.92
$ws.Range("E48").Value = 161.91999999999999
$ws.Range("F48").Value = 171.92

# Row 49
$ws.Range("A49").Value = 45961
$ws.Range("D49").Value = 166.31
$ws.Range("E49").Value = 161.08000000000001
$ws.Range("F49").Value = 171.08

# Row 50
$ws.Range("A50").Value = 45961
$ws.Range("D50").Value = 166.31
$ws.Range("E50").Value = 161.26
$ws.Range("F50").Value = 171.26

# Row 54
$ws.Range("A54").Value = 45962
$ws.Range("D54").Value = 184
$ws.Range("E54").Value = 176.13
$ws.Range("F54").Value = 186.13

# Row 55
$ws.Range("A55").Value = 45962
$ws.Range("D55").Value = 171.66

# Row 56
$ws.Range("A56").Value = 45962
$ws.Range("D56").Value = 174.04

# Row 57
$ws.Range("A57").Value = 45962
$ws.Range("D57").Value = 173.8

# Row 58
$ws.Range("A58").Value = 45962
$ws.Range("D58").Value = 169.71

# Row 59
$ws.Range("A59").Value = 45962
$ws.Range("D59").Value = 176.44
$ws.Range("E59").Value = 174.43

# Row 60
$ws.Range("A60").Value = 45961
$ws.Range("D60").Value = 183.88
$ws.Range("E60").Value = 176.17
$ws.Range("F60").Value = 186.17

# Row 61
$ws.Range("A61").Value = 45961
$ws.Range("D61").Value = 171.53
$ws.Range("E61").Value = 173.86
$ws.Range("F61").Value = 183.86

# Row 62
$ws.Range("A62").Value = 45961
$ws.Range("D62").Value = 173.91

# Row 63
$ws.Range("A63").Value = 45961
$ws.Range("D63").Value = 173.65
$ws.Range("E63").Value = 168.13

# Row 64
$ws.Range("A64").Value = 45961
$ws.Range("D64").Value = 169.56
$ws.Range("E64").Value = 164.18
$ws.Range("F64").Value = 174.18

# Row 65
$ws.Range("A65").Value = 45961
$ws.Range("D65").Value = 176.29
$ws.Range("E65").Value = 174.46
